$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header changes ---------------------------------------------
# "Complexity" header is replaced by two headers: "Time Complexity" (G1)
# and a brand-new "Space Complexity" column (H1).
$ws.Range("G1").Value = "Time Complexity"
$ws.Range("H1").Value = "Space Complexity"

# --- Row 2: existing "Trapping Rain Water" entry -------------------------
# The author swapped the contents that used to sit under "Solution" (E2)
# and "Complexity" (G2): E2 now shows the O(N) complexity text, G2 now
# shows the link to the C# solution. The pre-existing hyperlink objects
# stay anchored to the same cells (C2, E2) - only their displayed text
# changes.
$ws.Range("E2").Value = "O(N)"
$ws.Range("G2").Value = "https://github.com/Gershon-Tadepalli/DS-Algo/blob/master/DS-AlgoPractice/DS-AlgoPractice/TrapWater.cs"
# New "Space Complexity" value for this row.
$ws.Range("H2").Value = "O(1)"

# --- Row 3: brand-new "Best Time to buy&sell stock" entry ---------------
$ws.Range("A3").Value = "Arrays"
$ws.Range("B3").Value = "Best Time to buy&sell stock"
$ws.Range("C3").Value = "https://leetcode.com/problems/best-time-to-buy-and-sell-stock/"
$ws.Range("D3").Value = "find sum of subarray( max+= arr[i]-arr[i-1])" + [char]10 + "then find the largest subarray sum"
# (E3 / "Solution" is intentionally left blank for this entry.)
$ws.Range("F3").Value = "Easy but concept is sum of subarrays"
$ws.Range("G3").Value = "https://github.com/Gershon-Tadepalli/DS-Algo/blob/master/DS-AlgoPractice/DS-AlgoPractice/TrapWater.cs"
$ws.Range("H3").Value = "O(1)"

# D3 gets the same wrap-text style used by D2.
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)

# C3 becomes a real hyperlink (like C2 / E2) pointing at the LeetCode
# problem page, styled like the other hyperlink cells.
$ws.Hyperlinks.Add($ws.Range("C3"), "https://leetcode.com/problems/best-time-to-buy-and-sell-stock/")
$ws.Range("C3").Style = "Hyperlink"

# Row 3 is tall enough to show the two wrapped lines in D3.
$ws.Rows.Item(3).RowHeight = 29

# --- Column widths (best effort - widened to fit the new content) -------
$ws.Columns.Item(2).ColumnWidth = 22.76
$ws.Columns.Item(3).ColumnWidth = 54.92
$ws.Columns.Item(6).ColumnWidth = 30.92
$ws.Columns.Item(7).ColumnWidth = 13.76
$ws.Columns.Item(8).ColumnWidth = 14.59

# --- View state -----------------------------------------------------------
$ws.Range("E3").Select()
